$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1580.5883
$ws.Range("I28").Value = 1816.9231
$ws.Range("J28").Value = 812.5
$ws.Range("K28").Value = 1816.9231
$ws.Range("L28").Value = 812.5
$ws.Range("M28").Value = -1331.9231
$ws.Range("N28").Value = -1782.5

$ws.Range("H62").Value = 80562.234
$ws.Range("I62").Value = 126500.5
$ws.Range("J62").Value = 7061
$ws.Range("K62").Value = 126500.5
$ws.Range("L62").Value = 7061
$ws.Range("M62").Value = -125876.5
$ws.Range("N62").Value = -8309

$ws.Range("H63").Value = 32500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 32500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 32500
$ws.Range("N63").Value = -33748

$ws.Range("H65").Value = 80562.234
$ws.Range("I65").Value = 126500.5
$ws.Range("J65").Value = 7061
$ws.Range("K65").Value = 632502.5
$ws.Range("L65").Value = 35305
$ws.Range("M65").Value = -629382.5
$ws.Range("N65").Value = -41545

$ws.Range("H66").Value = 32500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 32500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 97500
$ws.Range("N66").Value = -103740

$ws.Range("H76").Value = 26194398
$ws.Range("I76").Value = 39289460
$ws.Range("J76").Value = 4270.5713
$ws.Range("K76").Value = 39289460
$ws.Range("L76").Value = 4270.5713
$ws.Range("M76").Value = -39289145
$ws.Range("N76").Value = -4900.5713

$ws.Range("H79").Value = 26194398
$ws.Range("I79").Value = 39289460
$ws.Range("J79").Value = 4270.5713
$ws.Range("K79").Value = 39289460
$ws.Range("L79").Value = 4270.5713
$ws.Range("M79").Value = -39288368
$ws.Range("N79").Value = -6454.5713

$ws.Range("H86").Value = 125006160
$ws.Range("I86").Value = 4826.5
$ws.Range("J86").Value = 250007500
$ws.Range("K86").Value = 4826.5
$ws.Range("L86").Value = 250007500
$ws.Range("M86").Value = -3703.5

$ws.Range("H88").Value = 1765631
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 1985084.9
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 1985084.9
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -1985896.9

$ws.Range("H89").Value = 125006160
$ws.Range("I89").Value = 4826.5
$ws.Range("J89").Value = 250007500
$ws.Range("K89").Value = 24132.5
$ws.Range("L89").Value = 1250037500
$ws.Range("M89").Value = -18516.5

$ws.Range("H91").Value = 1765631
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 1985084.9
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 1985084.9
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -1987892.9

$ws.Range("H92").Value = 1037.0358
$ws.Range("I92").Value = 1148.2727
$ws.Range("J92").Value = 629.1667
$ws.Range("K92").Value = 1148.2727
$ws.Range("L92").Value = 629.1667
$ws.Range("M92").Value = 99.72730000000001
$ws.Range("N92").Value = -3125.1667

$ws.Range("H98").Value = 1466.9231
$ws.Range("I98").Value = 1307
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1307
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 191
$ws.Range("N98").Value = -4996

$ws.Range("H122").Value = 1466.9231
$ws.Range("I122").Value = 1307
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3921
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1471
$ws.Range("N122").Value = -10900

$ws.Range("H129").Value = 2849842
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 4630806
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 13892418
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -13902418

$ws.Range("H132").Value = 1432528.8
$ws.Range("I132").Value = 6175
$ws.Range("J132").Value = 3334333.8
$ws.Range("K132").Value = 18525
$ws.Range("L132").Value = 10003001.4
$ws.Range("M132").Value = -15995
$ws.Range("N132").Value = -10008061.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10098.265
$ws.Range("I32").Value = 9736.679
$ws.Range("J32").Value = 11785.667
$ws.Range("K32").Value = 9736.679
$ws.Range("L32").Value = 11785.667
$ws.Range("M32").Value = -9449.679
$ws.Range("N32").Value = -12359.667

$ws.Range("H45").Value = 2737.111
$ws.Range("I45").Value = 2454.25
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 2454.25
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -2077.25
$ws.Range("N45").Value = -5754

$ws.Range("H63").Value = 1781
$ws.Range("I63").Value = 1968.3334
$ws.Range("J63").Value = 1500
$ws.Range("K63").Value = 1968.3334
$ws.Range("L63").Value = 1500
$ws.Range("M63").Value = -1282.3334
$ws.Range("N63").Value = -2872

$ws.Range("H66").Value = 1781
$ws.Range("I66").Value = 1968.3334
$ws.Range("J66").Value = 1500
$ws.Range("K66").Value = 9841.666999999999
$ws.Range("L66").Value = 7500
$ws.Range("M66").Value = -6409.666999999999
$ws.Range("N66").Value = -14364

$ws.Range("H97").Value = 460.76923
$ws.Range("I97").Value = 498.0909
$ws.Range("J97").Value = 255.5
$ws.Range("K97").Value = 498.0909
$ws.Range("L97").Value = 255.5
$ws.Range("M97").Value = -2.090899999999976

$ws.Range("H122").Value = 1562.3636
$ws.Range("I122").Value = 1562.3636
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4687.0908
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2237.0908
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2682.4285
$ws.Range("I132").Value = 2430.2273
$ws.Range("J132").Value = 3109.2307
$ws.Range("K132").Value = 7290.6819
$ws.Range("L132").Value = 9327.6921
$ws.Range("M132").Value = -4760.6819
$ws.Range("N132").Value = -14387.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 625.82355
$ws.Range("I94").Value = 339.6
$ws.Range("J94").Value = 1034.7142
$ws.Range("K94").Value = 339.6
$ws.Range("L94").Value = 1034.7142
$ws.Range("M94").Value = 111.4
$ws.Range("N94").Value = -1936.7142

$ws.Range("H134").Value = 98648.28999999999
$ws.Range("I134").Value = 5022.222
$ws.Range("J134").Value = 168867.83
$ws.Range("K134").Value = 15066.666
$ws.Range("L134").Value = 506603.49
$ws.Range("M134").Value = -12531.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4389.8486
$ws.Range("I31").Value = 1336.9445
$ws.Range("J31").Value = 8053.3335
$ws.Range("K31").Value = 1336.9445
$ws.Range("L31").Value = 8053.3335
$ws.Range("M31").Value = -1041.9445
$ws.Range("N31").Value = -8643.333500000001

$ws.Range("H34").Value = 4389.8486
$ws.Range("I34").Value = 1336.9445
$ws.Range("J34").Value = 8053.3335
$ws.Range("K34").Value = 1336.9445
$ws.Range("L34").Value = 8053.3335
$ws.Range("M34").Value = -1134.9445
$ws.Range("N34").Value = -8457.333500000001

$ws.Range("H62").Value = 10073.333
$ws.Range("I62").Value = 11514.286
$ws.Range("J62").Value = 8812.5
$ws.Range("K62").Value = 11514.286
$ws.Range("L62").Value = 8812.5
$ws.Range("M62").Value = -10890.286
$ws.Range("N62").Value = -10060.5

$ws.Range("H65").Value = 10073.333
$ws.Range("I65").Value = 11514.286
$ws.Range("J65").Value = 8812.5
$ws.Range("K65").Value = 57571.43
$ws.Range("L65").Value = 44062.5
$ws.Range("M65").Value = -54451.43
$ws.Range("N65").Value = -50302.5

$ws.Range("H70").Value = 21066.666
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 21066.666
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 21066.666
$ws.Range("N70").Value = -21696.666

$ws.Range("H73").Value = 21066.666
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 21066.666
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 21066.666
$ws.Range("N73").Value = -23250.666

$ws.Range("H132").Value = 4999.3335
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20058.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 999.6667
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 999.5
$ws.Range("K49").Value = 3000
$ws.Range("L49").Value = 2998.5
$ws.Range("M49").Value = -2844
$ws.Range("N49").Value = -3310.5

$ws.Range("H92").Value = 312
$ws.Range("I92").Value = 312
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 936
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 312
$ws.Range("N92").ClearContents()

$ws.Range("H131").Value = 65948.5
$ws.Range("I131").Value = 200366
$ws.Range("J131").Value = 4849.636
$ws.Range("K131").Value = 601098
$ws.Range("L131").Value = 14548.908
$ws.Range("M131").Value = -596058
$ws.Range("N131").Value = -24628.908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2194133.5
$ws.Range("I122").Value = 2194133.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6582400.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6579950.5

$ws.Range("H126").Value = 2586.1428
$ws.Range("I126").Value = 2801.2222
$ws.Range("J126").Value = 2199
$ws.Range("K126").Value = 8403.6666
$ws.Range("L126").Value = 6597
$ws.Range("M126").Value = -5933.6666
$ws.Range("N126").Value = -11537

$ws.Range("H132").Value = 2814.4482
$ws.Range("I132").Value = 1995.6111
$ws.Range("J132").Value = 4154.364
$ws.Range("K132").Value = 5986.8333
$ws.Range("L132").Value = 12463.092
$ws.Range("M132").Value = -3456.8333
$ws.Range("N132").Value = -17523.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5224

$ws.Range("H22").Value = 380.5
$ws.Range("I22").Value = 420.25
$ws.Range("J22").Value = 301
$ws.Range("K22").Value = 420.25
$ws.Range("L22").Value = 301
$ws.Range("M22").Value = -125.25
$ws.Range("N22").Value = -891

$ws.Range("H27").Value = 380.5
$ws.Range("I27").Value = 420.25
$ws.Range("J27").Value = 301
$ws.Range("K27").Value = 420.25
$ws.Range("L27").Value = 301
$ws.Range("M27").Value = -313.25
$ws.Range("N27").Value = -515

$ws.Range("H46").Value = 470.54544
$ws.Range("I46").Value = 445.25
$ws.Range("J46").Value = 485
$ws.Range("K46").Value = 466.66666
$ws.Range("L46").Value = 485
$ws.Range("M46").Value = -257.25
$ws.Range("N46").Value = -861

$ws.Range("H132").Value = 2446.4
$ws.Range("I132").Value = 2198.2666
$ws.Range("J132").Value = 3190.8
$ws.Range("K132").Value = 6594.7998
$ws.Range("L132").Value = 9572.400000000001
$ws.Range("M132").Value = -4064.7998
$ws.Range("N132").Value = -14632.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1465.2
$ws.Range("I122").Value = 1309.909
$ws.Range("J122").Value = 1655
$ws.Range("K122").Value = 3929.727
$ws.Range("L122").Value = 4965
$ws.Range("M122").Value = -1479.727
$ws.Range("N122").Value = -9865

$ws.Range("H132").Value = 1528.8334
$ws.Range("I132").Value = 1116
$ws.Range("J132").Value = 3097.6
$ws.Range("K132").Value = 3348
$ws.Range("L132").Value = 9292.799999999999
$ws.Range("M132").Value = -818
$ws.Range("N132").Value = -14352.8
